# Resize one column and one row of the deck's table (the "Table 3" graphic
# frame, on the slide that lists Movies / similar tabular content).
#
# The PowerPoint object model reports table geometry in points, while the
# underlying OOXML stores EMUs (1 pt = 12700 EMU), so the target EMU
# widths/heights are converted to points before being written through the
# Table.Columns / Table.Rows collections.
#
#   3rd column width : 2629080 EMU -> 2872800 EMU
#   2nd row   height : 3911760 EMU ->  369720 EMU
#
# (The source commit also shows the internal GUIDs of the <a:fld type=
# "slidenum"> page-number placeholders changing on every slide layout and
# the slide master. Those ids are opaque, randomly-minted identifiers with
# no corresponding property anywhere in the PowerPoint object model -- the
# displayed "<#>" text/type is unchanged -- so there is no COM call that can
# reproduce that part of the diff; it looks like incidental churn from
# whatever internal tool re-saved the deck rather than a user-visible edit.)

$p = $ppt.ActivePresentation

$emuPerPoint = 12700
$newColWidthEmu = 2872800
$newRowHeightEmu = 369720

# Find the shape that hosts the table, searching every slide so the script
# keeps working even if the table isn't on slide 3.
$tableShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $candidate = $slide.Shapes.Item($shi)
        if ($candidate.HasTable) {
            $tableShape = $candidate
        }
    }
}

$tbl = $tableShape.Table

$tbl.Columns.Item(3).Width = $newColWidthEmu / $emuPerPoint
$tbl.Rows.Item(2).Height = $newRowHeightEmu / $emuPerPoint

Write-Output "Column 3 width (pt): $($tbl.Columns.Item(3).Width)"
Write-Output "Row 2 height (pt): $($tbl.Rows.Item(2).Height)"
